# Update the NATMI LR-pair sheet with newly recomputed TPM-derived statistics.
# Ligand (G,H) values depend only on the "Sending cluster", Receptor (M,N)
# values depend only on the "Target cluster"; the specificity columns
# (I,J,O,P) are each value normalized against the sum across the three
# clusters, and the edge columns (Q,R,S,T) are simple products of the
# ligand/receptor value pairs (Q=G*M, R=H*N, S=I*O, T=J*P).
# This table was regenerated from updated source TPM data, so the literal
# values below reflect the new recomputed numbers for every affected row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per spreadsheet row, keyed by column letter.
$rowUpdates = @{
    2  = @{ G=3.037522333333333;  H=9.112567;  I=0.1153015356242242; J=0.1153015356242242;
            M=0.8194946666666666; N=2.458484;  O=0.1466535424263973; P=0.1466535424263973;
            Q=2.489233352047556;  R=22.403100168428; S=0.01690937864649592; T=0.01690937864649593 }
    3  = @{ G=3.037522333333333;  H=9.112567;  I=0.1153015356242242; J=0.1153015356242242;
            O=0.6452529427684778; P=0.6452529427684778;
            Q=10.95224240118333;  R=98.57018161065;  S=0.07439865516725513; T=0.07439865516725513 }
    4  = @{ G=3.037522333333333;  H=9.112567;  I=0.1153015356242242; J=0.1153015356242242;
            N=3.488457; O=0.2080935148051249; P=0.2080935148051249;
            Q=3.532088682124334;  R=31.788798139119; S=0.02399350181047314; T=0.02399350181047314 }
    5  = @{ I=0.325995654495798;  J=0.325995654495798;
            M=0.8194946666666666; N=2.458484;  O=0.1466535424263973; P=0.1466535424263973;
            Q=7.037887669060887;  R=63.34098902154799; S=0.04780841754742068; T=0.04780841754742068 }
    6  = @{ I=0.325995654495798;  J=0.325995654495798;
            O=0.6452529427684778; P=0.6452529427684778;
            S=0.2103496553931496; T=0.2103496553931496 }
    7  = @{ I=0.325995654495798;  J=0.325995654495798;
            N=3.488457; O=0.2080935148051249; P=0.2080935148051249;
            Q=9.986385310764332;  R=89.87746779687899; S=0.06783758155522773; T=0.06783758155522773 }
    8  = @{ I=0.5587028098799778; J=0.5587028098799777;
            M=0.8194946666666666; N=2.458484;  O=0.1466535424263973; P=0.1466535424263973;
            Q=12.06177923569422;  R=108.556013121248; S=0.08193574623248071; T=0.08193574623248071 }
    9  = @{ I=0.5587028098799778; J=0.5587028098799777;
            O=0.6452529427684778; P=0.6452529427684778;
            Q=53.06996624893333;  S=0.3605046322080731; T=0.3605046322080729 }
    10 = @{ I=0.5587028098799778; J=0.5587028098799777;
            N=3.488457; O=0.2080935148051249; P=0.2080935148051249;
            Q=17.11501811978934;  R=154.035163078104; S=0.1162624314394241; T=0.116262431439424 }
}

foreach ($rowNum in $rowUpdates.Keys) {
    $cols = $rowUpdates[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $ws.Range("$colLetter$rowNum").Value = $cols[$colLetter]
    }
}
